$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New forecasted consumption values (column A) and new timestamps (column B)
# for rows 2..97 (the forecast run shifted forward by one day: 2025-11-08 -> 2025-11-09,
# consumption values updated to reflect the retrained model run).
$aVals = @(5200, 5150, 5110, 5070, 5030, 5010, 4990, 4970, 4960, 4950, 4950, 4950, 4950, 4950, 4950, 4960, 4990, 5010, 5030, 5060, 5100, 5120, 5140, 5180, 5230, 5260, 5290, 5330, 5360, 5380, 5390, 5390, 5380, 5370, 5340, 5310, 5260, 5220, 5180, 5140, 5090, 5060, 5040, 5020, 5000, 5000, 5000, 5010, 5040, 5060, 5080, 5100, 5140, 5190, 5250, 5320, 5400, 5480, 5560, 5660, 5770, 5890, 6020, 6140, 6290, 6410, 6510, 6590, 6660, 6700, 6710, 6710, 6700, 6690, 6670, 6640, 6590, 6530, 6490, 6420, 6330, 6240, 6170, 6070, 5940, 5810, 5670, 5550, 5460, 5350, 5230, 5120, 5160, 5100, 5050, 5000)
$bVals = @(45970, 45970.01041666666, 45970.02083333334, 45970.03125, 45970.04166666666, 45970.05208333334, 45970.0625, 45970.07291666666, 45970.08333333334, 45970.09375, 45970.10416666666, 45970.11458333334, 45970.125, 45970.13541666666, 45970.14583333334, 45970.15625, 45970.16666666666, 45970.17708333334, 45970.1875, 45970.19791666666, 45970.20833333334, 45970.21875, 45970.22916666666, 45970.23958333334, 45970.25, 45970.26041666666, 45970.27083333334, 45970.28125, 45970.29166666666, 45970.30208333334, 45970.3125, 45970.32291666666, 45970.33333333334, 45970.34375, 45970.35416666666, 45970.36458333334, 45970.375, 45970.38541666666, 45970.39583333334, 45970.40625, 45970.41666666666, 45970.42708333334, 45970.4375, 45970.44791666666, 45970.45833333334, 45970.46875, 45970.47916666666, 45970.48958333334, 45970.5, 45970.51041666666, 45970.52083333334, 45970.53125, 45970.54166666666, 45970.55208333334, 45970.5625, 45970.57291666666, 45970.58333333334, 45970.59375, 45970.60416666666, 45970.61458333334, 45970.625, 45970.63541666666, 45970.64583333334, 45970.65625, 45970.66666666666, 45970.67708333334, 45970.6875, 45970.69791666666, 45970.70833333334, 45970.71875, 45970.72916666666, 45970.73958333334, 45970.75, 45970.76041666666, 45970.77083333334, 45970.78125, 45970.79166666666, 45970.80208333334, 45970.8125, 45970.82291666666, 45970.83333333334, 45970.84375, 45970.85416666666, 45970.86458333334, 45970.875, 45970.88541666666, 45970.89583333334, 45970.90625, 45970.91666666666, 45970.92708333334, 45970.9375, 45970.94791666666, 45970.95833333334, 45970.96875, 45970.97916666666, 45970.98958333334)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}
